# Applies the Dec-30-2023 "cryptos list" refresh: updated prices (col D) and
# 1h volume deltas (col E) for most rows, plus a rank swap between
# WEMIXToken and Stellar (rows 34/35, including coin name + link).
#
# Price cells (col D) hold values that look numeric ("0.970", "42.458.19",
# "1.613.13", ...) but must stay literal text (trailing zeros, multi-dot
# thousands groupings, etc. would be mangled/lost as a real number). Setting
# NumberFormat to "@" (Text) on those cells before writing the value forces
# Excel to keep the exact string instead of coercing it to a Double.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell($ref, $value) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $value
}

function Set-TextCell($ref, $value) {
    $ws.Range($ref).Value = $value
}

# row -> @{ D = new price; E = new volume-change text }
# (rows omitted here had no D/E change in this update)
$rowUpdates = [ordered]@{
    2  = @{ D = '42.458.19';  E = '  -1.30%  ' }
    3  = @{ D = '2.304.57';   E = '  -2.76%  ' }
    5  = @{ D = '319.98';     E = '  +1.18%  ' }
    6  = @{ D = '104.48';     E = '  -4.25%  ' }
    7  = @{ D = '0.632';      E = '  -0.99%  ' }
    8  = @{             E = '  +0.06%  ' }
    9  = @{ D = '0.614';      E = '  -0.77%  ' }
    10 = @{ D = '39.83';      E = '  -2.61%  ' }
    11 = @{ D = '0.0913';     E = '  -1.86%  ' }
    12 = @{ D = '8.38';       E = '  -2.65%  ' }
    13 = @{             E = '  -0.41%  ' }
    14 = @{ D = '0.970';      E = '  -3.89%  ' }
    15 = @{ D = '15.42';      E = '  -3.84%  ' }
    16 = @{ D = '2.652.63';   E = '  -3.02%  ' }
    17 = @{ D = '2.298.19';   E = '  -3.24%  ' }
    18 = @{ D = '42.397.30';  E = '  -1.51%  ' }
    19 = @{ D = '7.42';       E = '  -3.53%  ' }
    20 = @{             E = '  -0.76%  ' }
    21 = @{ D = '3.66';       E = '  +1.44%  ' }
    22 = @{ D = '73.66';      E = '  -3.87%  ' }
    23 = @{ D = '280.92';     E = '  +3.81%  ' }
    24 = @{ D = '10.96';      E = '  +14.20%  ' }
    25 = @{ D = '2.28';       E = '  -2.70%  ' }
    26 = @{             E = '  -0.19%  ' }
    27 = @{ D = '10.91';      E = '  -4.69%  ' }
    28 = @{ D = '2.38';       E = '  +5.49%  ' }
    29 = @{ D = '23.03';      E = '  -2.17%  ' }
    30 = @{ D = '36.42';      E = '  -1.23%  ' }
    31 = @{ D = '164.41';     E = '  -2.00%  ' }
    32 = @{ D = '0.0877';     E = '  -3.95%  ' }
    33 = @{ D = '5.90';       E = '  -3.85%  ' }
    # 34/35: WEMIXToken and Stellar swap rank; handled separately below.
    36 = @{ D = '0.114';      E = '  -6.16%  ' }
    37 = @{ D = '4.60';       E = '  -2.27%  ' }
    38 = @{ D = '0.0350';     E = '  -3.19%  ' }
    39 = @{ D = '3.78';       E = '  -1.55%  ' }
    40 = @{ D = '2.80';       E = '  +4.03%  ' }
    41 = @{ D = '100.21';     E = '  -4.95%  ' }
    42 = @{             E = '  -3.73%  ' }
    43 = @{ D = '69.61';      E = '  -2.83%  ' }
    44 = @{ D = '0.228';      E = '  -4.68%  ' }
    45 = @{             E = '  -0.06%  ' }
    46 = @{ D = '12.09';      E = '  -4.47%  ' }
    47 = @{ D = '112.34';     E = '  -2.05%  ' }
    48 = @{ D = '77.77';      E = '  -3.34%  ' }
    49 = @{ D = '8.96';       E = '  -2.21%  ' }
    50 = @{ D = '5.31';       E = '  -4.51%  ' }
    51 = @{ D = '1.613.13';   E = '  +1.86%  ' }
}

foreach ($row in $rowUpdates.Keys) {
    $update = $rowUpdates[$row]
    if ($update.ContainsKey('D')) {
        Set-PriceCell "D$row" $update.D
    }
    if ($update.ContainsKey('E')) {
        Set-TextCell "E$row" $update.E
    }
}

# Rows 34 & 35 swap places in the ranking: WEMIXToken drops below Stellar.
Set-TextCell  "B34" 'Stellar'
Set-TextCell  "C34" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-PriceCell "D34" '0.137'
Set-TextCell  "E34" '  +3.76%  '

Set-TextCell  "B35" 'WEMIXToken'
Set-TextCell  "C35" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-PriceCell "D35" '2.74'
Set-TextCell  "E35" '  -6.55%  '
